$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Step 1: back up the current last-row's formatting (border/numfmt pattern used
#     for the table's final row) onto a scratch row far away, so we can re-apply it
#     to the new final row (517) once the table has grown. ---
$ws.Range("A501:O501").Copy()
$ws.Range("A1000:O1000").PasteSpecial(-4122) | Out-Null

# --- Step 2: grow the table by 16 rows (501 -> 517) using the table object model,
#     which is the supported way to extend a ListObject's range. ---
for ($i = 0; $i -lt 16; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# --- Step 3: re-establish the alternating row-banding format. Row 501 (formerly the
#     last row) becomes a normal "odd" banded row; the new rows 502-516 alternate
#     even/odd using rows 500 and 499 as formatting templates; row 517 becomes the
#     new last row and gets the special closing-border format saved in step 1. ---
$ws.Range("A499:O499").Copy()
$ws.Range("A501:O501").PasteSpecial(-4122) | Out-Null

for ($r = 502; $r -le 516; $r++) {
    if ($r % 2 -eq 0) {
        $ws.Range("A500:O500").Copy()
    } else {
        $ws.Range("A499:O499").Copy()
    }
    $ws.Range("A" + $r + ":O" + $r).PasteSpecial(-4122) | Out-Null
}

$ws.Range("A1000:O1000").Copy()
$ws.Range("A517:O517").PasteSpecial(-4122) | Out-Null

# clean up the scratch row used to stash the closing-border format
$ws.Range("A1000:O1000").ClearContents() | Out-Null
$ws.Range("A1000:O1000").ClearFormats() | Out-Null

# --- Step 4: write the new response data (rows 502-517). Row 501's own values are
#     unchanged by this edit (only its banding/border style moved), so it is left
#     alone here. ---
    # Row 502
    $ws.Range("A502").Value = 45599.00320989583
    $ws.Range("B502").Value = "h20202564@glab.hallym.ac.kr"
    $ws.Range("C502").Value = "디지털미디어콘텐츠"
    $ws.Range("D502").Value = 20202564.0
    $ws.Range("E502").Value = "이호연"
    $ws.Range("F502").Value = "랜덤화"
    $ws.Range("G502").Value = "28 vs 71"
    $ws.Range("H502").Value = "NFIP 설계의 대조군 집단"
    $ws.Range("I502").Value = "Red"
    $ws.Range("J502").Value = "가. 10센트"
    $ws.Range("K502").Value = "나. 5분"
    $ws.Range("L502").Value = "나. 47일"

    # Row 503
    $ws.Range("A503").Value = 45599.00509670139
    $ws.Range("B503").Value = "you72460601@gmail.com"
    $ws.Range("C503").Value = "디지털인문예술"
    $ws.Range("D503").Value = 20201721.0
    $ws.Range("E503").Value = "유지원"
    $ws.Range("F503").Value = "가짜약 대조군"
    $ws.Range("G503").Value = "28 vs 71"
    $ws.Range("H503").Value = "NFIP 설계의 백신 접종 집단"
    $ws.Range("I503").Value = "Red"
    $ws.Range("J503").Value = "가. 10센트"
    $ws.Range("K503").Value = "나. 5분"
    $ws.Range("L503").Value = "가. 24일"

    # Row 504
    $ws.Range("A504").Value = 45599.01188731482
    $ws.Range("B504").Value = "kimeunji0512@naver.com"
    $ws.Range("C504").Value = "언어청각학부"
    $ws.Range("D504").Value = 20243910.0
    $ws.Range("E504").Value = "김은지"
    $ws.Range("F504").Value = "랜덤화"
    $ws.Range("G504").Value = "28 vs 25"
    $ws.Range("H504").Value = "NFIP 설계의 백신 접종 집단"
    $ws.Range("I504").Value = "Black"
    $ws.Range("M504").Value = "가. 5센트"
    $ws.Range("N504").Value = "가. 5분"
    $ws.Range("O504").Value = "가. 47일"

    # Row 505
    $ws.Range("A505").Value = 45599.03350423611
    $ws.Range("B505").Value = "jytoto33@naver.com"
    $ws.Range("C505").Value = "언어청각학부"
    $ws.Range("D505").Value = 20243912.0
    $ws.Range("E505").Value = "김지윤"
    $ws.Range("F505").Value = "랜덤화"
    $ws.Range("G505").Value = "28 vs 71"
    $ws.Range("H505").Value = "NFIP 설계의 대조군 집단"
    $ws.Range("I505").Value = "Red"
    $ws.Range("J505").Value = "나. 5센트"
    $ws.Range("K505").Value = "나. 5분"
    $ws.Range("L505").Value = "나. 47일"

    # Row 506
    $ws.Range("A506").Value = 45599.04710865741
    $ws.Range("B506").Value = "hanseoyun392@gmail.com"
    $ws.Range("C506").Value = "언어청각학부 "
    $ws.Range("D506").Value = 20243964.0
    $ws.Range("E506").Value = "한서윤"
    $ws.Range("F506").Value = "랜덤화"
    $ws.Range("G506").Value = "28 vs 71"
    $ws.Range("H506").Value = "NFIP 설계의 백신 접종 집단"
    $ws.Range("I506").Value = "Black"
    $ws.Range("M506").Value = "가. 5센트"
    $ws.Range("N506").Value = "가. 5분"
    $ws.Range("O506").Value = "가. 47일"

    # Row 507
    $ws.Range("A507").Value = 45599.05038975694
    $ws.Range("B507").Value = "benjamin27@naver.com"
    $ws.Range("C507").Value = "디지털미디어 콘텐츠"
    $ws.Range("D507").Value = 20212583.0
    $ws.Range("E507").Value = "최재혁"
    $ws.Range("F507").Value = "가짜약 대조군"
    $ws.Range("G507").Value = "28 vs 71"
    $ws.Range("H507").Value = "NFIP 설계의 대조군 집단"
    $ws.Range("I507").Value = "Red"
    $ws.Range("J507").Value = "나. 5센트"
    $ws.Range("K507").Value = "나. 5분"
    $ws.Range("L507").Value = "나. 47일"

    # Row 508
    $ws.Range("A508").Value = 45599.051567766204
    $ws.Range("B508").Value = "hsjenny99@gmail.com"
    $ws.Range("C508").Value = "소프트웨어학부"
    $ws.Range("D508").Value = 20245246.0
    $ws.Range("E508").Value = "전소현"
    $ws.Range("F508").Value = "랜덤화"
    $ws.Range("G508").Value = "25 vs 54"
    $ws.Range("H508").Value = "플라시보 컨트롤 설계의 생리식염수 접종 집단"
    $ws.Range("I508").Value = "Red"
    $ws.Range("J508").Value = "나. 5센트"
    $ws.Range("K508").Value = "나. 5분"
    $ws.Range("L508").Value = "나. 47일"

    # Row 509
    $ws.Range("A509").Value = 45599.085377615746
    $ws.Range("B509").Value = "moon050123@naver.com"
    $ws.Range("C509").Value = "중국학과"
    $ws.Range("D509").Value = 20241518.0
    $ws.Range("E509").Value = "문서원"
    $ws.Range("F509").Value = "랜덤화"
    $ws.Range("G509").Value = "28 vs 71"
    $ws.Range("H509").Value = "NFIP 설계의 대조군 집단"
    $ws.Range("I509").Value = "Red"
    $ws.Range("J509").Value = "가. 10센트"
    $ws.Range("K509").Value = "나. 5분"
    $ws.Range("L509").Value = "가. 24일"

    # Row 510
    $ws.Range("A510").Value = 45599.09575731482
    $ws.Range("B510").Value = "ehdus1113kim@naver.com"
    $ws.Range("C510").Value = "법학과"
    $ws.Range("D510").Value = 20232705.0
    $ws.Range("E510").Value = "김도연"
    $ws.Range("F510").Value = "이중눈가림"
    $ws.Range("G510").Value = "28 vs 25"
    $ws.Range("H510").Value = "플라시보 컨트롤 설계의 생리식염수 접종 집단"
    $ws.Range("I510").Value = "Red"
    $ws.Range("J510").Value = "나. 5센트"
    $ws.Range("K510").Value = "나. 5분"
    $ws.Range("L510").Value = "나. 47일"

    # Row 511
    $ws.Range("A511").Value = 45599.09997949074
    $ws.Range("B511").Value = "vldzmgha0609@naver.com"
    $ws.Range("C511").Value = "디지털미디어콘텐츠"
    $ws.Range("D511").Value = 20227039.0
    $ws.Range("E511").Value = "김지수"
    $ws.Range("F511").Value = "이중눈가림"
    $ws.Range("G511").Value = "28 vs 71"
    $ws.Range("H511").Value = "플라시보 컨트롤 설계의 백신 접종 집단"
    $ws.Range("I511").Value = "Black"
    $ws.Range("M511").Value = "나. 10센트"
    $ws.Range("N511").Value = "나. 100분"
    $ws.Range("O511").Value = "나. 24일"

    # Row 512
    $ws.Range("A512").Value = 45599.11320365741
    $ws.Range("B512").Value = "a35142191@gmail.com"
    $ws.Range("C512").Value = "데이터사이언스학부"
    $ws.Range("D512").Value = 20243241.0
    $ws.Range("E512").Value = "이윤재"
    $ws.Range("F512").Value = "랜덤화"
    $ws.Range("G512").Value = "28 vs 71"
    $ws.Range("H512").Value = "NFIP 설계의 대조군 집단"
    $ws.Range("I512").Value = "Black"
    $ws.Range("M512").Value = "나. 10센트"
    $ws.Range("N512").Value = "가. 5분"
    $ws.Range("O512").Value = "나. 24일"

    # Row 513
    $ws.Range("A513").Value = 45599.14028357639
    $ws.Range("B513").Value = "krcar1002@gmail.com"
    $ws.Range("C513").Value = "심리학과"
    $ws.Range("D513").Value = 20242109.0
    $ws.Range("E513").Value = "김재호"
    $ws.Range("F513").Value = "랜덤화"
    $ws.Range("G513").Value = "28 vs 71"
    $ws.Range("H513").Value = "NFIP 설계의 대조군 집단"
    $ws.Range("I513").Value = "Black"
    $ws.Range("M513").Value = "나. 10센트"
    $ws.Range("N513").Value = "나. 100분"
    $ws.Range("O513").Value = "나. 24일"

    # Row 514
    $ws.Range("A514").Value = 45599.14806766204
    $ws.Range("B514").Value = "hyeonyonga@naver.com"
    $ws.Range("C514").Value = "화학과"
    $ws.Range("D514").Value = 20243415.0
    $ws.Range("E514").Value = "안현용"
    $ws.Range("F514").Value = "랜덤화"
    $ws.Range("G514").Value = "28 vs 71"
    $ws.Range("H514").Value = "NFIP 설계의 백신 접종 집단"
    $ws.Range("I514").Value = "Red"
    $ws.Range("J514").Value = "나. 5센트"
    $ws.Range("K514").Value = "나. 5분"
    $ws.Range("L514").Value = "나. 47일"

    # Row 515
    $ws.Range("A515").Value = 45599.226665532406
    $ws.Range("B515").Value = "jud050207@gmail.com"
    $ws.Range("C515").Value = "소프트웨어학부"
    $ws.Range("D515").Value = 20245252.0
    $ws.Range("E515").Value = "정의돈"
    $ws.Range("F515").Value = "랜덤화"
    $ws.Range("G515").Value = "28 vs 71"
    $ws.Range("H515").Value = "NFIP 설계의 대조군 집단"
    $ws.Range("I515").Value = "Red"
    $ws.Range("J515").Value = "가. 10센트"
    $ws.Range("K515").Value = "나. 5분"
    $ws.Range("L515").Value = "나. 47일"

    # Row 516
    $ws.Range("A516").Value = 45599.22714829861
    $ws.Range("B516").Value = "jyj111212@naver.com"
    $ws.Range("C516").Value = "인문학부"
    $ws.Range("D516").Value = 20241083.0
    $ws.Range("E516").Value = "장예지"
    $ws.Range("F516").Value = "가짜약 대조군"
    $ws.Range("G516").Value = "28 vs 46"
    $ws.Range("H516").Value = "NFIP 설계의 백신 접종 집단"
    $ws.Range("I516").Value = "Black"
    $ws.Range("M516").Value = "나. 10센트"
    $ws.Range("N516").Value = "가. 5분"
    $ws.Range("O516").Value = "가. 47일"

    # Row 517
    $ws.Range("A517").Value = 45599.23857908565
    $ws.Range("B517").Value = "dmsdn6462@naver.com"
    $ws.Range("C517").Value = "심리학과"
    $ws.Range("D517").Value = 20202106.0
    $ws.Range("E517").Value = "김은우"
    $ws.Range("F517").Value = "가짜약 대조군"
    $ws.Range("G517").Value = "28 vs 71"
    $ws.Range("H517").Value = "NFIP 설계의 대조군 집단"
    $ws.Range("I517").Value = "Black"
    $ws.Range("M517").Value = "가. 5센트"
    $ws.Range("N517").Value = "가. 5분"
    $ws.Range("O517").Value = "가. 47일"

# --- Step 5: each response only answers the Red-track (4_R/5_R/6_R -> J:L) or the
#     Black-track (4_B/5_B/6_B -> M:O) questions depending on the "Red or Black"
#     column (I); clear out the unused track's cells so they stay empty, matching
#     the pattern used throughout the rest of the table. ---
    $ws.Range("M502:O502").ClearContents() | Out-Null
    $ws.Range("M503:O503").ClearContents() | Out-Null
    $ws.Range("J504:L504").ClearContents() | Out-Null
    $ws.Range("M505:O505").ClearContents() | Out-Null
    $ws.Range("J506:L506").ClearContents() | Out-Null
    $ws.Range("M507:O507").ClearContents() | Out-Null
    $ws.Range("M508:O508").ClearContents() | Out-Null
    $ws.Range("M509:O509").ClearContents() | Out-Null
    $ws.Range("M510:O510").ClearContents() | Out-Null
    $ws.Range("J511:L511").ClearContents() | Out-Null
    $ws.Range("J512:L512").ClearContents() | Out-Null
    $ws.Range("J513:L513").ClearContents() | Out-Null
    $ws.Range("M514:O514").ClearContents() | Out-Null
    $ws.Range("M515:O515").ClearContents() | Out-Null
    $ws.Range("J516:L516").ClearContents() | Out-Null
    $ws.Range("J517:L517").ClearContents() | Out-Null